# EZ-4133 imail templates update again
#
# 1. The lone "_GoBack" bookmark that currently sits in front of the
#    "Notice of event of default" heading is moved: it now sits inside the
#    signature block, between the literal "(" and "formerly known as ..."
#    text that gets appended after "EZBob Ltd". Adding a bookmark with the
#    same name elsewhere automatically relocates it (Word only allows one
#    bookmark per name), so no separate delete step is required.
# 2. Three new runs are appended to the "EZBob Ltd" signature paragraph so
#    it reads: "EZBob Ltd. (formerly known as Orange Money Ltd.)" - split
#    across the same run boundaries as the authored edit.

$d = $word.ActiveDocument

# --- locate the "EZBob" / " Ltd" signature paragraph -----------------
# "EZBob Ltd" also appears in the footer address block and in the
# recipient clause near the top of the letter, so anchor the search right
# after the "Yours sincerely" valediction that immediately precedes the
# signature block we need.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Yours sincerely", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$sigSearch = $d.Range($anchor.End, $d.Content.End)
$sigSearch.Find.Execute("EZBob Ltd", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPos = $sigSearch.End

# --- append the three new runs, matching the "EZBob"/" Ltd" formatting ---
$run1 = $d.Range($insertPos, $insertPos)
$run1.InsertAfter(". ")
$run1.Font.Name = "Open Sans Light"
$run1.Font.Color = 6710886

$pos2 = $run1.End
$run2 = $d.Range($pos2, $pos2)
$run2.InsertAfter("(")
$run2.Font.Name = "Open Sans Light"
$run2.Font.Color = 6710886

$pos3 = $run2.End
$run3 = $d.Range($pos3, $pos3)
$run3.InsertAfter("formerly known as Orange Money Ltd.)")
$run3.Font.Name = "Open Sans Light"
$run3.Font.Color = 6710886

# --- move the "_GoBack" bookmark to sit right after the new "(" ------
# Re-find the freshly inserted text rather than reusing the InsertAfter
# ranges directly above: adding a zero-length bookmark exactly at the
# stale end-anchor of the run that was just typed into can silently no-op,
# so re-anchor via a fresh Find instead.
$bmSearch = $d.Content.Duplicate
$bmSearch.Find.Execute("(formerly known", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $bmSearch.Start + 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
